# Auto-generated Excel COM-interop script
# Applies market-price / leve-profit data refresh across all 8 profession sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (59 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H15").Value = 1978.9773
$ws.Range("I15").Value = 1978.9773
$ws.Range("K15").Value = 5936.9319
$ws.Range("M15").Value = -5767.9319
$ws.Range("H17").Value = 1112731.6
$ws.Range("J17").Value = 1112731.6
$ws.Range("L17").Value = 3338194.8
$ws.Range("N17").Value = -3338530.8
$ws.Range("H18").Value = 417.83334
$ws.Range("I18").Value = 392.52173
$ws.Range("K18").Value = 392.52173
$ws.Range("M18").Value = -108.52173
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 1662.7273
$ws.Range("J32").Value = 1712
$ws.Range("L32").Value = 1712
$ws.Range("N32").Value = -2364
$ws.Range("H53").Value = 622.0952
$ws.Range("I53").Value = 434.2
$ws.Range("K53").Value = 434.2
$ws.Range("M53").Value = 202.8
$ws.Range("H98").Value = 3366061.8
$ws.Range("I98").Value = 3499111.2
$ws.Range("K98").Value = 3499111.2
$ws.Range("M98").Value = -3497613.2
$ws.Range("H103").Value = 735.8333
$ws.Range("I103").Value = 736.7778
$ws.Range("K103").Value = 2210.3334
$ws.Range("M103").Value = -1624.3334
$ws.Range("H107").Value = 20837544
$ws.Range("I107").Value = 12504005
$ws.Range("K107").Value = 12504005
$ws.Range("M107").Value = -12502085
$ws.Range("H122").Value = 3366061.8
$ws.Range("I122").Value = 3499111.2
$ws.Range("K122").Value = 10497333.6
$ws.Range("M122").Value = -10494883.6
$ws.Range("H132").Value = 1592.3695
$ws.Range("J132").Value = 2733.75
$ws.Range("L132").Value = 8201.25
$ws.Range("N132").Value = -13261.25
$ws.Range("H137").Value = 28728.953
$ws.Range("I137").Value = 32439.473
$ws.Range("J137").Value = 6465.8335
$ws.Range("K137").Value = 97318.41900000001
$ws.Range("L137").Value = 19397.5005
$ws.Range("M137").Value = -94768.41900000001
$ws.Range("N137").Value = -24497.5005
$ws.Range("H138").Value = 3586.3677
$ws.Range("I138").Value = 2056.5881
$ws.Range("K138").Value = 6169.7643
$ws.Range("M138").Value = -1029.7643

# --- Sheet: ARM (66 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3385
$ws.Range("I2").Value = 2600
$ws.Range("J2").Value = 4955
$ws.Range("K2").Value = 2600
$ws.Range("L2").Value = 4955
$ws.Range("M2").Value = -2487
$ws.Range("N2").Value = -5181
$ws.Range("H32").Value = 6553.6665
$ws.Range("I32").Value = 4787.4507
$ws.Range("K32").Value = 4787.4507
$ws.Range("M32").Value = -4500.4507
$ws.Range("H45").Value = 1979.2
$ws.Range("I45").Value = 1979.2
$ws.Range("K45").Value = 1979.2
$ws.Range("M45").Value = -1602.2
$ws.Range("H61").Value = 4038.074
$ws.Range("I61").Value = 1966.15
$ws.Range("K61").Value = 1966.15
$ws.Range("M61").Value = -1754.15
$ws.Range("H74").Value = 59159.145
$ws.Range("I74").Value = 67954.5
$ws.Range("J74").Value = 6387
$ws.Range("K74").Value = 67954.5
$ws.Range("L74").Value = 6387
$ws.Range("M74").Value = -67080.5
$ws.Range("N74").Value = -8135
$ws.Range("H77").Value = 59159.145
$ws.Range("I77").Value = 67954.5
$ws.Range("J77").Value = 6387
$ws.Range("K77").Value = 339772.5
$ws.Range("L77").Value = 31935
$ws.Range("M77").Value = -335404.5
$ws.Range("N77").Value = -40671
$ws.Range("H97").Value = 1138.375
$ws.Range("I97").Value = 484.66666
$ws.Range("J97").Value = 2227.889
$ws.Range("K97").Value = 484.66666
$ws.Range("L97").Value = 2227.889
$ws.Range("M97").Value = 11.33334000000002
$ws.Range("N97").Value = -3219.889
$ws.Range("H103").Value = 10000
$ws.Range("J103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("N103").Value = -12344
$ws.Range("H116").Value = 3385
$ws.Range("I116").Value = 2600
$ws.Range("J116").Value = 4955
$ws.Range("K116").Value = 2600
$ws.Range("L116").Value = 4955
$ws.Range("M116").Value = -306
$ws.Range("N116").Value = -9543
$ws.Range("H121").Value = 189975
$ws.Range("J121").Value = 189975
$ws.Range("L121").Value = 189975
$ws.Range("N121").Value = -193469
$ws.Range("H122").Value = 2209.5
$ws.Range("I122").Value = 2218.08
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 6654.24
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -4204.24
$ws.Range("N122").Value = -10885
$ws.Range("H136").Value = 4038.074
$ws.Range("I136").Value = 1966.15
$ws.Range("K136").Value = 5898.450000000001
$ws.Range("M136").Value = -3348.450000000001

# --- Sheet: BSM (33 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3385
$ws.Range("I3").Value = 2600
$ws.Range("J3").Value = 4955
$ws.Range("K3").Value = 2600
$ws.Range("L3").Value = 4955
$ws.Range("M3").Value = -2486
$ws.Range("N3").Value = -5183
$ws.Range("H22").Value = 834.6
$ws.Range("I22").Value = 834.6
$ws.Range("K22").Value = 834.6
$ws.Range("M22").Value = -661.6
$ws.Range("H99").Value = 2988.2903
$ws.Range("I99").Value = 3393.52
$ws.Range("J99").Value = 1299.8334
$ws.Range("K99").Value = 3393.52
$ws.Range("L99").Value = 1299.8334
$ws.Range("M99").Value = -1895.52
$ws.Range("N99").Value = -4295.8334
$ws.Range("H107").Value = 1786
$ws.Range("I107").Value = 1997
$ws.Range("J107").Value = 1575
$ws.Range("K107").Value = 1997
$ws.Range("L107").Value = 1575
$ws.Range("M107").Value = -77
$ws.Range("N107").Value = -5415
$ws.Range("H117").Value = 97999.664
$ws.Range("J117").Value = 97999.664
$ws.Range("L117").Value = 97999.664
$ws.Range("N117").Value = -107177.664
$ws.Range("H141").Value = 104995
$ws.Range("J141").Value = 104995
$ws.Range("L141").Value = 104995
$ws.Range("N141").Value = -115355

# --- Sheet: CRP (84 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 228.28572
$ws.Range("I7").Value = 207.25
$ws.Range("K7").Value = 207.25
$ws.Range("M7").Value = -94.25
$ws.Range("H16").Value = 2810.3635
$ws.Range("I16").Value = 2351
$ws.Range("K16").Value = 2351
$ws.Range("M16").Value = -2064
$ws.Range("H31").Value = 246934.34
$ws.Range("I31").Value = 305014.28
$ws.Range("J31").Value = 7354.625
$ws.Range("K31").Value = 305014.28
$ws.Range("L31").Value = 7354.625
$ws.Range("M31").Value = -304719.28
$ws.Range("N31").Value = -7944.625
$ws.Range("H34").Value = 246934.34
$ws.Range("I34").Value = 305014.28
$ws.Range("J34").Value = 7354.625
$ws.Range("K34").Value = 305014.28
$ws.Range("L34").Value = 7354.625
$ws.Range("M34").Value = -304812.28
$ws.Range("N34").Value = -7758.625
$ws.Range("H58").Value = 1875.5483
$ws.Range("I58").Value = 1724.9642
$ws.Range("J58").Value = 3281
$ws.Range("K58").Value = 1724.9642
$ws.Range("L58").Value = 3281
$ws.Range("M58").Value = -1521.9642
$ws.Range("N58").Value = -3687
$ws.Range("H86").Value = 4188.5
$ws.Range("I86").Value = 3055.5
$ws.Range("J86").Value = 4868.3
$ws.Range("K86").Value = 3055.5
$ws.Range("L86").Value = 4868.3
$ws.Range("M86").Value = -1932.5
$ws.Range("N86").Value = -7114.3
$ws.Range("H89").Value = 4188.5
$ws.Range("I89").Value = 3055.5
$ws.Range("J89").Value = 4868.3
$ws.Range("K89").Value = 15277.5
$ws.Range("L89").Value = 24341.5
$ws.Range("M89").Value = -9661.5
$ws.Range("N89").Value = -35573.5
$ws.Range("H105").Value = 5916.8335
$ws.Range("I105").Value = 2138.2727
$ws.Range("J105").Value = 8104.421
$ws.Range("K105").Value = 2138.2727
$ws.Range("L105").Value = 8104.421
$ws.Range("M105").Value = -391.2727
$ws.Range("N105").Value = -11598.421
$ws.Range("H107").Value = 7268.524
$ws.Range("I107").Value = 2746.6667
$ws.Range("K107").Value = 2746.6667
$ws.Range("M107").Value = -826.6667000000002
$ws.Range("H113").Value = 2810.3635
$ws.Range("I113").Value = 2351
$ws.Range("K113").Value = 2351
$ws.Range("M113").Value = -181
$ws.Range("H121").Value = 48323
$ws.Range("J121").Value = 48326
$ws.Range("L121").Value = 48326
$ws.Range("N121").Value = -50946
$ws.Range("H122").Value = 2808.5833
$ws.Range("I122").Value = 2154.818
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6464.454000000001
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -4014.454000000001
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 3714.353
$ws.Range("I132").Value = 2199.1667
$ws.Range("K132").Value = 6597.500100000001
$ws.Range("M132").Value = -4067.500100000001
$ws.Range("H134").Value = 7211.421
$ws.Range("I134").Value = 8559.267
$ws.Range("K134").Value = 25677.801
$ws.Range("M134").Value = -23142.801
$ws.Range("H136").Value = 1875.5483
$ws.Range("I136").Value = 1724.9642
$ws.Range("J136").Value = 3281
$ws.Range("K136").Value = 5174.892599999999
$ws.Range("L136").Value = 9843
$ws.Range("M136").Value = -2624.892599999999
$ws.Range("N136").Value = -14943

# --- Sheet: CUL (53 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H12").Value = 3435.7273
$ws.Range("I12").Value = 442.75
$ws.Range("K12").Value = 1328.25
$ws.Range("M12").Value = -1155.25
$ws.Range("H50").Value = 799.75
$ws.Range("I50").Value = 649.5
$ws.Range("K50").Value = 1948.5
$ws.Range("M50").Value = -1467.5
$ws.Range("H53").Value = 799.75
$ws.Range("I53").Value = 649.5
$ws.Range("K53").Value = 1948.5
$ws.Range("M53").Value = -1467.5
$ws.Range("H68").Value = 9261413
$ws.Range("I68").Value = 16667782
$ws.Range("K68").Value = 50003346
$ws.Range("M68").Value = -50002535
$ws.Range("H71").Value = 9261413
$ws.Range("I71").Value = 16667782
$ws.Range("K71").Value = 150010038
$ws.Range("M71").Value = -150005982
$ws.Range("H86").Value = 870.4286
$ws.Range("I86").Value = 832.1667
$ws.Range("K86").Value = 2496.5001
$ws.Range("M86").Value = -1310.5001
$ws.Range("H89").Value = 870.4286
$ws.Range("I89").Value = 832.1667
$ws.Range("K89").Value = 7489.5003
$ws.Range("M89").Value = -1561.5003
$ws.Range("H122").Value = 328.42856
$ws.Range("I122").Value = 314.5
$ws.Range("J122").Value = 334
$ws.Range("K122").Value = 2830.5
$ws.Range("L122").Value = 3006
$ws.Range("M122").Value = -380.5
$ws.Range("N122").Value = -7906
$ws.Range("H131").Value = 7354316
$ws.Range("I131").Value = 62500796
$ws.Range("K131").Value = 187502388
$ws.Range("M131").Value = -187497348
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

# --- Sheet: GSM (58 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29039.666
$ws.Range("I15").Value = 23999.5
$ws.Range("J15").Value = 39120
$ws.Range("K15").Value = 23999.5
$ws.Range("L15").Value = 39120
$ws.Range("M15").Value = -23711.5
$ws.Range("N15").Value = -39696
$ws.Range("H81").Value = 29039.666
$ws.Range("I81").Value = 23999.5
$ws.Range("J81").Value = 39120
$ws.Range("K81").Value = 23999.5
$ws.Range("L81").Value = 39120
$ws.Range("M81").Value = -23001.5
$ws.Range("N81").Value = -41116
$ws.Range("H84").Value = 29039.666
$ws.Range("I84").Value = 23999.5
$ws.Range("J84").Value = 39120
$ws.Range("K84").Value = 71998.5
$ws.Range("L84").Value = 117360
$ws.Range("M84").Value = -67006.5
$ws.Range("N84").Value = -127344
$ws.Range("H107").Value = 92754
$ws.Range("I107").Value = 101809.3
$ws.Range("K107").Value = 101809.3
$ws.Range("M107").Value = -99889.3
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 2081.2
$ws.Range("I113").Value = 1955
$ws.Range("J113").Value = 2165.3333
$ws.Range("K113").Value = 1955
$ws.Range("L113").Value = 2165.3333
$ws.Range("M113").Value = 215
$ws.Range("N113").Value = -6505.3333
$ws.Range("H114").Value = 112333
$ws.Range("J114").Value = 112333
$ws.Range("L114").Value = 112333
$ws.Range("N114").Value = -121011
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H126").Value = 8936
$ws.Range("I126").Value = 12549.5
$ws.Range("J126").Value = 4118
$ws.Range("K126").Value = 37648.5
$ws.Range("L126").Value = 12354
$ws.Range("M126").Value = -35178.5
$ws.Range("N126").Value = -17294
$ws.Range("H132").Value = 33214.64
$ws.Range("I132").Value = 43146.445
$ws.Range("J132").Value = 3419.2222
$ws.Range("K132").Value = 129439.335
$ws.Range("L132").Value = 10257.6666
$ws.Range("M132").Value = -126909.335
$ws.Range("N132").Value = -15317.6666

# --- Sheet: LTW (29 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23494.223
$ws.Range("I40").Value = 23494.223
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 23494.223
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -23358.223
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 2158.484
$ws.Range("I46").Value = 2551.111
$ws.Range("K46").Value = 2551.111
$ws.Range("M46").Value = -2363.111
$ws.Range("H55").Value = 677.1053000000001
$ws.Range("I55").Value = 805.86664
$ws.Range("K55").Value = 805.86664
$ws.Range("M55").Value = -632.86664
$ws.Range("H132").Value = 3943.4375
$ws.Range("I132").Value = 2662.1614
$ws.Range("J132").Value = 5147.0605
$ws.Range("K132").Value = 7986.4842
$ws.Range("L132").Value = 15441.1815
$ws.Range("M132").Value = -5456.4842
$ws.Range("N132").Value = -20501.1815
$ws.Range("H136").Value = 1964.1562
$ws.Range("I136").Value = 1553.381
$ws.Range("J136").Value = 2748.3635
$ws.Range("K136").Value = 4660.143
$ws.Range("L136").Value = 8245.0905
$ws.Range("M136").Value = -2110.143
$ws.Range("N136").Value = -13345.0905

# --- Sheet: WVR (23 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1630.7273
$ws.Range("I100").Value = 1641.5555
$ws.Range("K100").Value = 3283.111
$ws.Range("M100").Value = -2742.111
$ws.Range("H107").Value = 1336.4546
$ws.Range("I107").Value = 1214.5555
$ws.Range("K107").Value = 3643.6665
$ws.Range("M107").Value = -1723.6665
$ws.Range("H122").Value = 2251.0908
$ws.Range("J122").Value = 3417.6
$ws.Range("L122").Value = 10252.8
$ws.Range("N122").Value = -15152.8
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 2054.5945
$ws.Range("I132").Value = 1572.2354
$ws.Range("J132").Value = 2464.6
$ws.Range("K132").Value = 4716.706200000001
$ws.Range("L132").Value = 7393.799999999999
$ws.Range("M132").Value = -2186.706200000001
$ws.Range("N132").Value = -12453.8
